# The edit cyclically rotates the four data rows (2-5) of the "Artfynd"
# sheet: row 2's record moves down to become row 5, and rows 3,4,5 each
# shift up by one (3->2, 4->3, 5->4). All cell values/types are preserved
# verbatim for each record; only the row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that ever hold data in rows 2-5, split by native storage type so
# values round-trip through COM without Excel "helpfully" reinterpreting a
# numeric-looking or date-looking string as a real number/date.
$textCols = @("C","D","F","G","H","I","J","K","M","P","T","U","V","W","Y","Z","AA","AB","AI","AJ","AK","AO","AT","AW","AX","AY")
$numCols  = @("A","B","E","Q","R","S")
$boolCols = @("AD","AE","AG")

$firstDataRow = 2
$lastDataRow = 5

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($col in $textCols) {
        $snap[$col] = $ws.Range($col + $row).Value2
    }
    foreach ($col in $numCols) {
        $snap[$col] = $ws.Range($col + $row).Value2
    }
    foreach ($col in $boolCols) {
        $snap[$col] = $ws.Range($col + $row).Value2
    }
    return $snap
}

function Set-RowFromSnapshot($row, $snap) {
    foreach ($col in $textCols) {
        $cell = $ws.Range($col + $row)
        $val = $snap[$col]
        if ($val -eq $null) {
            # Source cell was entirely absent (not merely blank) - mirror
            # that exactly instead of leaving/forcing a "" string behind.
            $cur = $cell.Value2
            if (-not ($cur -eq $null)) {
                $cell.ClearContents()
            }
        } elseif ($val -eq "") {
            # Source cell existed but held an empty string. Only touch the
            # destination if it isn't already in that same state.
            $cur = $cell.Value2
            if ($cur -eq $null) {
                $cell.NumberFormat = "@"
                $cell.Value = $val
                $cell.Style = "Normal"
            }
        } else {
            # Force text storage so numeric-looking ("2","9") or
            # date-looking ("2018-04-29","00:00") strings are not
            # auto-converted to a number/date serial by Excel.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        }
    }
    foreach ($col in $numCols) {
        $ws.Range($col + $row).Value = $snap[$col]
    }
    foreach ($col in $boolCols) {
        $ws.Range($col + $row).Value = $snap[$col]
    }
}

# Capture the row-2 record before it gets overwritten.
$row2Snapshot = Get-RowSnapshot $firstDataRow

# Shift rows 3,4,5 up into 2,3,4.
for ($row = $firstDataRow; $row -lt $lastDataRow; $row++) {
    $sourceSnapshot = Get-RowSnapshot ($row + 1)
    Set-RowFromSnapshot $row $sourceSnapshot
}

# The original row-2 record becomes the new last row.
Set-RowFromSnapshot $lastDataRow $row2Snapshot
